# Add Q4-2022 data to the workbook:
#  1. Update the "总计" (summary) sheet: insert a new leading data row for
#     2022-Q4 and shift the existing quarters down by one row.
#  2. Insert a brand-new "2022-Q4" worksheet (cloned from "2022-Q3" so it
#     keeps identical headers/styling) positioned right after "总计", and
#     fill in the Q4 fund-holding figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet — push rows 2-5 down to 3-6 and write the new Q4 row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the new row 6 (previously empty) the same "A" column style as its
# neighbours before writing values into it.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q2"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 3.12

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 2.2

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 2.86

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 3.01

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 2.6

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet, cloned from "2022-Q3" and placed right before it.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Columns D:G hold text-formatted numbers in this workbook (e.g. "50.81"),
# so force the "@" text format before assigning, otherwise Excel would
# auto-coerce the numeric-looking strings into real numbers.
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "50.81"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "94.73"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "5.12"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "2.6015"
$q4.Range("H2").Value = 10

# Restore the original active tab ("2021-Q2" is selected in the source
# workbook) since copying a sheet makes the copy the active one.
$wb.Worksheets.Item("2021-Q2").Activate()

Write-Host "applied 2022-Q4 update"
